$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.832.43"
$ws.Range("E2").Value = "  -1.38%  "
$ws.Range("D3").Value = "3.370.01"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'572.52"
$ws.Range("E5").Value = "  -0.80%  "
$ws.Range("D6").Value = "'136.73"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "3.366.77"
$ws.Range("E8").Value = "  -0.67%  "
$ws.Range("E9").Value = "  -1.03%  "
$ws.Range("D10").Value = "'7.66"
$ws.Range("E10").Value = "  +2.23%  "
$ws.Range("E11").Value = "  -2.40%  "
$ws.Range("E12").Value = "  -1.77%  "
$ws.Range("D13").Value = "3.947.37"
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("D15").Value = "'25.87"
$ws.Range("E15").Value = "  +2.00%  "
$ws.Range("D16").Value = "3.371.01"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").Value = "'0.0000171"
$ws.Range("E17").Value = "  -3.19%  "
$ws.Range("D18").Value = "60.987.80"
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("D19").Value = "'13.91"
$ws.Range("E19").Value = "  -1.83%  "
$ws.Range("E20").Value = "  -0.87%  "
$ws.Range("D21").Value = "'9.38"
$ws.Range("E21").Value = "  -0.94%  "
$ws.Range("D22").Value = "'373.65"
$ws.Range("E22").Value = "  -3.51%  "
$ws.Range("D23").Value = "'0.551"
$ws.Range("D24").Value = "3.511.95"
$ws.Range("E24").Value = "  -0.66%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("E26").Value = "  -1.70%  "
$ws.Range("D27").Value = "'71.04"
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("D28").Value = "'1.77"
$ws.Range("E28").Value = "  +11.98%  "
$ws.Range("D29").Value = "'0.177"
$ws.Range("E29").Value = "  +9.64%  "
$ws.Range("D30").Value = "'7.45"
$ws.Range("E30").Value = "  -2.86%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").Value = "'8.10"
$ws.Range("E32").Value = "  -2.09%  "
$ws.Range("E33").Value = "  -0.89%  "
$ws.Range("D35").Value = "'23.64"
$ws.Range("E35").Value = "  +0.75%  "
$ws.Range("D36").Value = "'5.18"
$ws.Range("E36").Value = "  -4.20%  "
$ws.Range("D37").Value = "'6.86"
$ws.Range("E37").Value = "  -1.33%  "
$ws.Range("E38").Value = "  -1.07%  "
$ws.Range("D39").Value = "'164.49"
$ws.Range("E39").Value = "  +1.03%  "
$ws.Range("E40").Value = "  -3.06%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").Value = "'0.775"
$ws.Range("E42").Value = "  -0.99%  "
$ws.Range("D43").Value = "'41.53"
$ws.Range("E43").Value = "  -0.54%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "'1.70"
$ws.Range("E44").Value = "  -4.46%  "
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").Value = "'4.40"
$ws.Range("E45").Value = "  -1.03%  "
$ws.Range("E46").Value = "  -2.80%  "
$ws.Range("D47").Value = "'24.05"
$ws.Range("E47").Value = "  -2.22%  "
$ws.Range("D48").Value = "2.454.65"
$ws.Range("E48").Value = "  +4.27%  "
$ws.Range("D49").Value = "'6.79"
$ws.Range("E49").Value = "  -2.14%  "
$ws.Range("E50").Value = "  -1.69%  "
$ws.Range("D51").Value = "'2.41"
$ws.Range("E51").Value = "  +4.02%  "
